$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-DiffFormula($row) {
    $f = '=IF(B' + $row + '="","-",IF(C' + $row + '="","-",B' + $row + '-C' + $row + '))'
    $ws.Range("D$row").Formula = $f
}

# ---------------------------------------------------------------------
# First, write the new "Place" labels (column A) in the exact order the
# new strings were first introduced in the authored workbook, so that
# the shared-string table grows in the matching order (154..161):
#   154 Letter Get, 155 Both Boss HP = 30, 156 Both Boss HP = 20,
#   157 Both Boss HP = 10, 158 Both Boss HP = 0, 159 Both Letter Get,
#   160 Start run, 161 Boss 3 appears
# ---------------------------------------------------------------------
$ws.Range("A187").Value = "Letter Get"
$ws.Range("A188").Value = "Both Boss HP = 30"
$ws.Range("A189").Value = "Both Boss HP = 20"
$ws.Range("A190").Value = "Both Boss HP = 10"
$ws.Range("A193").Value = "Both Boss HP = 0"
$ws.Range("A194").Value = "Both Letter Get"
$ws.Range("A191").Value = "Start run"
$ws.Range("A192").Value = "Start run"
$ws.Range("A195").Value = "Boss 3 appears"

# --- Row 187: Letter Get ---
$ws.Range("B187").Value = 161921
$ws.Range("C187").Value = 150349
Set-DiffFormula 187

# --- Row 188: Both Boss HP = 30 ---
$ws.Range("B188").Value = 162219
$ws.Range("C188").Value = 150691
Set-DiffFormula 188

# --- Row 189: Both Boss HP = 20 ---
$ws.Range("B189").Value = 162499
$ws.Range("C189").Value = 150931
Set-DiffFormula 189

# --- Row 190: Both Boss HP = 10 (no C value) ---
$ws.Range("B190").Value = 162755
$ws.Range("C190").ClearContents()
Set-DiffFormula 190

# --- Row 191: Start run (only A + C, no B, no D) ---
$ws.Range("B191").ClearContents()
$ws.Range("C191").Value = 151267
$ws.Range("D191").ClearContents()

# --- Row 192: Start run (only A + C, no B, no D) ---
$ws.Range("B192").ClearContents()
$ws.Range("C192").Value = 151435
$ws.Range("D192").ClearContents()

# --- Row 193: Both Boss HP = 0 ---
$ws.Range("B193").Value = 163167
$ws.Range("C193").Value = 151418
Set-DiffFormula 193

# --- Row 194: Both Letter Get (no C value) ---
$ws.Range("B194").Value = 163502
$ws.Range("C194").ClearContents()
Set-DiffFormula 194

# --- Row 195: Boss 3 appears ---
$ws.Range("B195").Value = 163623
$ws.Range("C195").Value = 151874
Set-DiffFormula 195

# --- Rows 196-197 stay empty (nothing to do) ---

# --- Row 198: only D formula, no A/B/C ---
Set-DiffFormula 198

# --- Row 199: end level (re-uses existing shared string) ---
$ws.Range("A199").Value = "end level"
$ws.Range("B199").Value = 175028
Set-DiffFormula 199

# --- Row 200: boss fight end (white screen) (re-uses existing shared string) ---
$ws.Range("A200").Value = "boss fight end (white screen)"
$ws.Range("B200").Value = 179257
Set-DiffFormula 200

# --- Update view: scroll position & selection to match authored state ---
$ws.Range("C196").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 180
